# Applies the "added test for symmetric layers fm44" change:
#  - H47/I47 flip from "in progress" to "complete"
#  - a new row 49 is appended describing the fm44 test case
#  - sheet view selection moves to F49 (topLeftCell scrolls to A22)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# fm42's status (H47/I47) moves from "in progress" to "complete"
$ws.Range("H47").Value = "complete"
$ws.Range("I47").Value = "complete"

# New row for the fm44 symmetric-layers test case. Copy row 48's formatting
# down to row 49 first (format-only paste) so the new row picks up the same
# cell styles (s="5" text cells, s="7" right-aligned numeric cell), then
# overwrite the values.
$ws.Range("B48:I48").Copy()
$ws.Range("B49:I49").PasteSpecial(-4122)   # xlPasteFormats

# The description is written before the "fm44" label so the two new shared
# strings land in the same order as the reference workbook
# (description=148, fm44=149).
$ws.Range("C49").Value = "Multiple accounts with same number of layers (policies) per account. Account level output"
$ws.Range("B49").Value = "fm44"
$ws.Range("D49").Value = "All"
$ws.Range("E49").Value = "2,12,14"
$ws.Range("F49").Value = 3
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = "complete"
$ws.Range("I49").Value = "complete"

# Update the selection/scroll position to reflect the new active cell
$ws.Range("F49").Select()
$excel.ActiveWindow.ScrollRow = 22
